# Update the LR-pairs data table (Fgf7-Nrp1) with refreshed TPM-derived values.
# The sending-cluster set grows from {ECs, FAPs, MuSCs} to {ECs, FAPs, Inflammatory-Mac, MuSCs},
# so the table grows from 15 data rows (16 incl. header) to 20 data rows (21 incl. header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row, in column order A..T:
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value, Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$rows = @(
    @('ECs','Fgf7','Nrp1','ECs',1,0.3333333333333333,0.1658776666666667,0.497633,0.01379803068909966,0.01379803068909966,3,1,86.89540866666668,260.686226,0.319779657009892,0.3197796570098919,14.41400763367311,129.726068703058,0.004412329521172252,0.004412329521172252),
    @('ECs','Fgf7','Nrp1','FAPs',1,0.3333333333333333,0.1658776666666667,0.497633,0.01379803068909966,0.01379803068909966,3,1,54.34621166666667,163.038635,0.1999969065479545,0.1999969065479545,9.014822783439445,81.133405050955,0.002759563454273673,0.002759563454273673),
    @('ECs','Fgf7','Nrp1','Inflammatory-Mac',1,0.3333333333333333,0.1658776666666667,0.497633,0.01379803068909966,0.01379803068909966,3,1,60.92601633333334,182.778049,0.224210932487692,0.224210932487692,10.10626542866856,90.95638885801701,0.003093669327296825,0.003093669327296825),
    @('ECs','Fgf7','Nrp1','MuSCs',1,0.3333333333333333,0.1658776666666667,0.497633,0.01379803068909966,0.01379803068909966,3,1,7.809668333333332,23.429005,0.02873998867505581,0.02873998867505581,1.295449560573889,11.659046045165,0.0003965552457427967,0.0003965552457427968),
    @('ECs','Fgf7','Nrp1','Resolving-Mac',1,0.3333333333333333,0.1658776666666667,0.497633,0.01379803068909966,0.01379803068909966,3,1,61.75795633333333,185.273869,0.2272725152794058,0.2272725152794058,10.24426569467522,92.19839125207699,0.003135913140614113,0.003135913140614113),
    @('FAPs','Fgf7','Nrp1','ECs',3,1,11.51839566666667,34.555187,0.9581228147923823,0.9581228147923824,3,1,86.89540866666668,260.686226,0.319779657009892,0.3197796570098919,1000.895698639363,9008.061287754264,0.3063881850876602,0.3063881850876602),
    @('FAPs','Fgf7','Nrp1','FAPs',3,1,11.51839566666667,34.555187,0.9581228147923823,0.9581228147923824,3,1,54.34621166666667,163.038635,0.1999969065479545,0.1999969065479545,625.9811689610829,5633.830520649745,0.1916215990514952,0.1916215990514952),
    @('FAPs','Fgf7','Nrp1','Inflammatory-Mac',3,1,11.51839566666667,34.555187,0.9581228147923823,0.9581228147923824,3,1,60.92601633333334,182.778049,0.224210932487692,0.224210932487692,701.7699625211294,6315.929662690164,0.2148216097423322,0.2148216097423322),
    @('FAPs','Fgf7','Nrp1','MuSCs',3,1,11.51839566666667,34.555187,0.9581228147923823,0.9581228147923824,3,1,7.809668333333332,23.429005,0.02873998867505581,0.02873998867505581,89.95484988877055,809.593648998935,0.02753643884644567,0.02753643884644567),
    @('FAPs','Fgf7','Nrp1','Resolving-Mac',3,1,11.51839566666667,34.555187,0.9581228147923823,0.9581228147923824,3,1,61.75795633333333,185.273869,0.2272725152794058,0.2272725152794058,711.3525766120559,6402.173189508503,0.217754982064449,0.217754982064449),
    @('Inflammatory-Mac','Fgf7','Nrp1','ECs',2,0.6666666666666666,0.1594243333333333,0.478273,0.01326122972505393,0.01326122972505393,3,1,86.89540866666668,260.686226,0.319779657009892,0.3197796570098919,13.85324259641089,124.679183367698,0.004240671493007129,0.004240671493007128),
    @('Inflammatory-Mac','Fgf7','Nrp1','FAPs',2,0.6666666666666666,0.1594243333333333,0.478273,0.01326122972505393,0.01326122972505393,3,1,54.34621166666667,163.038635,0.1999969065479545,0.1999969065479545,8.664108564150556,77.976977077355,0.002652204922032567,0.002652204922032567),
    @('Inflammatory-Mac','Fgf7','Nrp1','Inflammatory-Mac',2,0.6666666666666666,0.1594243333333333,0.478273,0.01326122972505393,0.01326122972505393,3,1,60.92601633333334,182.778049,0.224210932487692,0.224210932487692,9.713089536597446,87.417805829377,0.00297331268258784,0.00297331268258784),
    @('Inflammatory-Mac','Fgf7','Nrp1','MuSCs',2,0.6666666666666666,0.1594243333333333,0.478273,0.01326122972505393,0.01326122972505393,3,1,7.809668333333332,23.429005,0.02873998867505581,0.02873998867505581,1.245051167596111,11.205460508365,0.0003811275921153633,0.0003811275921153634),
    @('Inflammatory-Mac','Fgf7','Nrp1','Resolving-Mac',2,0.6666666666666666,0.1594243333333333,0.478273,0.01326122972505393,0.01326122972505393,3,1,61.75795633333333,185.273869,0.2272725152794058,0.2272725152794058,9.845721016470778,92.39341398700803,0.003013913035311029,0.003013913035311029),
    @('MuSCs','Fgf7','Nrp1','ECs',2,0.6666666666666666,0.1781386666666667,0.534416,0.01481792479346402,0.01481792479346403,3,1,86.89540866666668,260.686226,0.319779657009892,0.3197796570098919,15.47943223933511,139.314890154016,0.0047384709080523,0.0047384709080523),
    @('MuSCs','Fgf7','Nrp1','FAPs',2,0.6666666666666666,0.1781386666666667,0.534416,0.01481792479346402,0.01481792479346403,3,1,54.34621166666667,163.038635,0.1999969065479545,0.1999969065479545,9.681161684684445,87.13045516216,0.002963539120153043,0.002963539120153043),
    @('MuSCs','Fgf7','Nrp1','Inflammatory-Mac',2,0.6666666666666666,0.1781386666666667,0.534416,0.01481792479346402,0.01481792479346403,3,1,60.92601633333334,182.778049,0.224210932487692,0.224210932487692,10.85327931493156,97.67951383438401,0.003322340735475059,0.003322340735475059),
    @('MuSCs','Fgf7','Nrp1','MuSCs',2,0.6666666666666666,0.1781386666666667,0.534416,0.01481792479346402,0.01481792479346403,3,1,7.809668333333332,23.429005,0.02873998867505581,0.02873998867505581,1.391203904008889,12.52083513608,0.0004258669907519848,0.0004258669907519849),
    @('MuSCs','Fgf7','Nrp1','Resolving-Mac',2,0.6666666666666666,0.1781386666666667,0.534416,0.01481792479346402,0.01481792479346403,3,1,61.75795633333333,185.273869,0.2272725152794058,0.2272725152794058,11.00147999727822,99.01331997550399,0.003367707039031639,0.003367707039031639)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowValues = $rows[$i]
    $excelRow = $startRow + $i
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($excelRow, $col).Value = $rowValues[$col - 1]
    }
}
